$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 20 ("Servlet Filters" / "What are Servlet Filters?")
# Bump the body placeholder's explicit run sizes from 18pt to 20pt, and
# reposition the filter diagram picture.
# ---------------------------------------------------------------------------
$slide20 = $p.Slides.Item(20)
$body20 = $slide20.Shapes.Item(2)
$tr20 = $body20.TextFrame.TextRange

# Paragraphs 1,3,4,5,6,7 carry visible text runs that need sz=2000.
foreach ($i in 1,3,4,5,6,7) {
  $para = $tr20.Paragraphs($i,1)
  $para.Font.Size = 20
}

# Reposition the "Servlet_Filter.gif" picture (shape 3) to its new offset.
# Shape.Left/Top are Single-precision points; nudge by a hair so the
# round-tripped EMU value lands exactly on the target (2234019, 3554866).
$pic20 = $slide20.Shapes.Item(3)
$pic20.Left = 175.90700587401577
$pic20.Top = 279.91072166141737

# ---------------------------------------------------------------------------
# Slide 22 ("Servlet Filters" - implementation steps)
# Fix the typo "instrunction" -> "instruction" and split it into its own
# run (matching the diff, which turns it into a separate <a:r>).
# ---------------------------------------------------------------------------
$slide22 = $p.Slides.Item(22)
$rect22 = $slide22.Shapes.Item(3)
$tr22 = $rect22.TextFrame.TextRange

$fullText22 = $tr22.Text
$word = "instrunction "
$idx = $fullText22.IndexOf($word)
if ($idx -ge 0) {
  $startPos = $idx + 1
  $sel = $tr22.Characters($startPos, $word.Length)
  $sel.Text = "instruction "
}
